$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: gpt_answer changes from "a" to "d"; gpt_correct flips from TRUE to FALSE
$ws.Range("H3").Value = "d"
$ws.Range("I3").Value = $false

# Row 45: gpt_answer changes from "b" to "a"; gpt_correct flips from TRUE to FALSE
$ws.Range("H45").Value = "a"
$ws.Range("I45").Value = $false

# Row 57: gpt_answer changes from "a" to "c"; gpt_correct flips from FALSE to TRUE
$ws.Range("H57").Value = "c"
$ws.Range("I57").Value = $true
